$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 280
$ws1.Range("F4").Value = 317
$ws1.Range("F6").Value = 654
$ws1.Range("F7").Value = 268
$ws1.Range("F11").Value = 132
$ws1.Range("F12").Value = 3356
$ws1.Range("F13").Value = 104
$ws1.Range("F14").Value = 74
$ws1.Range("F16").Value = 35
$ws1.Range("F17").Value = 49
$ws1.Range("F18").Value = 569
$ws1.Range("F19").Value = 40
$ws1.Range("F20").Value = 667
$ws1.Range("F21").Value = 196
$ws1.Range("F22").Value = 108
$ws1.Range("F24").Value = 46
$ws1.Range("F26").Value = 2372
$ws1.Range("F27").Value = 4915
$ws1.Range("F30").Value = 472
$ws1.Range("F31").Value = 1260
$ws1.Range("F32").Value = 267
$ws1.Range("F33").Value = 2182
$ws1.Range("F35").Value = 482
$ws1.Range("F37").Value = 73
$ws1.Range("F38").Value = 153
$ws1.Range("F39").Value = 304
$ws1.Range("F40").Value = 450
$ws1.Range("F41").Value = 769
$ws1.Range("F42").Value = 22
$ws1.Range("F44").Value = 27
$ws1.Range("F45").Value = 452

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 64

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 280
$ws4.Range("F4").Value = 317
$ws4.Range("F6").Value = 654
$ws4.Range("F7").Value = 268
$ws4.Range("F11").Value = 132
$ws4.Range("F12").Value = 3356
$ws4.Range("F13").Value = 104
$ws4.Range("F14").Value = 74
$ws4.Range("F16").Value = 64
$ws4.Range("F17").Value = 35
$ws4.Range("F18").Value = 49
$ws4.Range("F19").Value = 569
$ws4.Range("F20").Value = 40
$ws4.Range("F21").Value = 667
$ws4.Range("F22").Value = 196
$ws4.Range("F23").Value = 108
$ws4.Range("F25").Value = 46
$ws4.Range("F27").Value = 2372
$ws4.Range("F28").Value = 4915
$ws4.Range("F31").Value = 472
$ws4.Range("F32").Value = 1260
$ws4.Range("F33").Value = 267
$ws4.Range("F34").Value = 2182
$ws4.Range("F36").Value = 482
$ws4.Range("F38").Value = 73
$ws4.Range("F39").Value = 153
$ws4.Range("F40").Value = 304
$ws4.Range("F41").Value = 450
$ws4.Range("F42").Value = 769
$ws4.Range("F43").Value = 22
$ws4.Range("F45").Value = 27
$ws4.Range("F46").Value = 452
